# Auto-generated Excel COM-interop script
# Applies scheduled market-data value updates across multiple sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 2461.923  # H40: 2545.9092 -> 2461.923
$ws.Cells.Item(40, 10).Value = 2750.625  # J40: 2857.8572 -> 2750.625
$ws.Cells.Item(40, 12).Value = 2750.625  # L40: 2857.8572 -> 2750.625
$ws.Cells.Item(40, 14).Value = -3100.625  # N40: -3207.8572 -> -3100.625
$ws.Cells.Item(51, 8).Value = 62399.39  # H51: 43999.383 -> 62399.39
$ws.Cells.Item(51, 9).Value = 10199  # I51: 8324.375 -> 10199
$ws.Cells.Item(51, 10).Value = 88499.586  # J51: 59854.945 -> 88499.586
$ws.Cells.Item(51, 11).Value = 10199  # K51: 8324.375 -> 10199
$ws.Cells.Item(51, 12).Value = 88499.586  # L51: 59854.945 -> 88499.586
$ws.Cells.Item(51, 13).Value = -9715  # M51: -7840.375 -> -9715
$ws.Cells.Item(51, 14).Value = -89467.586  # N51: -60822.945 -> -89467.586
$ws.Cells.Item(55, 8).Value = 307.22223  # H55: 298.42105 -> 307.22223
$ws.Cells.Item(55, 9).Value = 266.47058  # I55: 259.44446 -> 266.47058
$ws.Cells.Item(55, 11).Value = 266.47058  # K55: 259.44446 -> 266.47058
$ws.Cells.Item(55, 13).Value = -52.47057999999998  # M55: -45.44445999999999 -> -52.47057999999998
$ws.Cells.Item(135, 8).Value = 1465.4286  # H135: 1471.4286 -> 1465.4286
$ws.Cells.Item(135, 9).Value = 1205.5454  # I135: 1276.1 -> 1205.5454
$ws.Cells.Item(135, 10).Value = 2418.3333  # J135: 1959.75 -> 2418.3333
$ws.Cells.Item(135, 11).Value = 10849.9086  # K135: 11484.9 -> 10849.9086
$ws.Cells.Item(135, 12).Value = 21764.9997  # L135: 17637.75 -> 21764.9997
$ws.Cells.Item(135, 13).Value = -8314.908599999999  # M135: -8949.9 -> -8314.908599999999
$ws.Cells.Item(135, 14).Value = -26834.9997  # N135: -22707.75 -> -26834.9997
$ws.Cells.Item(138, 8).Value = 2349.2122  # H138: 3011.423 -> 2349.2122
$ws.Cells.Item(138, 9).Value = 1784.6666  # I138: 1784.3529 -> 1784.6666
$ws.Cells.Item(138, 10).Value = 2560.9167  # J138: 5329.222 -> 2560.9167
$ws.Cells.Item(138, 11).Value = 5353.9998  # K138: 5353.0587 -> 5353.9998
$ws.Cells.Item(138, 12).Value = 7682.750100000001  # L138: 15987.666 -> 7682.750100000001
$ws.Cells.Item(138, 13).Value = -213.9997999999996  # M138: -213.0587000000005 -> -213.9997999999996
$ws.Cells.Item(138, 14).Value = -17962.7501  # N138: -26267.666 -> -17962.7501
# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 3079.8  # H45: 2899.8333 -> 3079.8
$ws.Cells.Item(45, 9).Value = 1499.5  # I45: 1666.3334 -> 1499.5
$ws.Cells.Item(45, 11).Value = 1499.5  # K45: 1666.3334 -> 1499.5
$ws.Cells.Item(45, 13).Value = -1122.5  # M45: -1289.3334 -> -1122.5
$ws.Cells.Item(122, 8).Value = 16850.264  # H122: 15257.381 -> 16850.264
$ws.Cells.Item(122, 9).Value = 12572.25  # I122: 11189.223 -> 12572.25
$ws.Cells.Item(122, 11).Value = 37716.75  # K122: 33567.669 -> 37716.75
$ws.Cells.Item(122, 13).Value = -35266.75  # M122: -31117.669 -> -35266.75
# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(12, 8).Value = 1298.5555  # H12: 1143.5454 -> 1298.5555
$ws.Cells.Item(12, 9).Value = 363.33334  # I12: 471.5 -> 363.33334
$ws.Cells.Item(12, 10).Value = 1766.1666  # J12: 1527.5714 -> 1766.1666
$ws.Cells.Item(12, 11).Value = 363.33334  # K12: 471.5 -> 363.33334
$ws.Cells.Item(12, 12).Value = 1766.1666  # L12: 1527.5714 -> 1766.1666
$ws.Cells.Item(12, 13).Value = -195.33334  # M12: -303.5 -> -195.33334
$ws.Cells.Item(12, 14).Value = -2102.1666  # N12: -1863.5714 -> -2102.1666
$ws.Cells.Item(64, 8).Value = 1753.25  # H64: 1669 -> 1753.25
$ws.Cells.Item(64, 10).Value = 2006.5  # J64: 2007 -> 2006.5
$ws.Cells.Item(64, 12).Value = 2006.5  # L64: 2007 -> 2006.5
$ws.Cells.Item(64, 14).Value = -2456.5  # N64: -2457 -> -2456.5
$ws.Cells.Item(67, 8).Value = 1753.25  # H67: 1669 -> 1753.25
$ws.Cells.Item(67, 10).Value = 2006.5  # J67: 2007 -> 2006.5
$ws.Cells.Item(67, 12).Value = 2006.5  # L67: 2007 -> 2006.5
$ws.Cells.Item(67, 14).Value = -3566.5  # N67: -3567 -> -3566.5
$ws.Cells.Item(94, 8).Value = 3438.2222  # H94: 2994.6 -> 3438.2222
$ws.Cells.Item(94, 9).Value = 1983.3334  # I94: 1789.8 -> 1983.3334
$ws.Cells.Item(94, 10).Value = 4165.6665  # J94: 4199.4 -> 4165.6665
$ws.Cells.Item(94, 11).Value = 1983.3334  # K94: 1789.8 -> 1983.3334
$ws.Cells.Item(94, 12).Value = 4165.6665  # L94: 4199.4 -> 4165.6665
$ws.Cells.Item(94, 13).Value = -1532.3334  # M94: -1338.8 -> -1532.3334
$ws.Cells.Item(94, 14).Value = -5067.6665  # N94: -5101.4 -> -5067.6665
$ws.Cells.Item(107, 8).Value = 0  # H107: 1125 -> 0
$ws.Cells.Item(107, 9).Value = 0  # I107: 1125 -> 0
$ws.Cells.Item(107, 11).Value = 0  # K107: 1125 -> 0
$ws.Cells.Item(107, 13).ClearContents()  # M107: was 795
$ws.Cells.Item(135, 8).Value = 72332  # H135: 72498.5 -> 72332
$ws.Cells.Item(135, 10).Value = 72332  # J135: 72498.5 -> 72332
$ws.Cells.Item(135, 12).Value = 72332  # L135: 72498.5 -> 72332
$ws.Cells.Item(135, 14).Value = -82472  # N135: -82638.5 -> -82472
# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3336.366  # H31: 3338.4634 -> 3336.366
$ws.Cells.Item(31, 9).Value = 2322.2  # I31: 2374.2068 -> 2322.2
$ws.Cells.Item(31, 10).Value = 6102.273  # J31: 5668.75 -> 6102.273
$ws.Cells.Item(31, 11).Value = 2322.2  # K31: 2374.2068 -> 2322.2
$ws.Cells.Item(31, 12).Value = 6102.273  # L31: 5668.75 -> 6102.273
$ws.Cells.Item(31, 13).Value = -2027.2  # M31: -2079.2068 -> -2027.2
$ws.Cells.Item(31, 14).Value = -6692.273  # N31: -6258.75 -> -6692.273
$ws.Cells.Item(34, 8).Value = 3336.366  # H34: 3338.4634 -> 3336.366
$ws.Cells.Item(34, 9).Value = 2322.2  # I34: 2374.2068 -> 2322.2
$ws.Cells.Item(34, 10).Value = 6102.273  # J34: 5668.75 -> 6102.273
$ws.Cells.Item(34, 11).Value = 2322.2  # K34: 2374.2068 -> 2322.2
$ws.Cells.Item(34, 12).Value = 6102.273  # L34: 5668.75 -> 6102.273
$ws.Cells.Item(34, 13).Value = -2120.2  # M34: -2172.2068 -> -2120.2
$ws.Cells.Item(34, 14).Value = -6506.273  # N34: -6072.75 -> -6506.273
$ws.Cells.Item(99, 8).Value = 10799.517  # H99: 10818.936 -> 10799.517
$ws.Cells.Item(99, 10).Value = 14098.875  # J99: 14136.5 -> 14098.875
$ws.Cells.Item(99, 12).Value = 14098.875  # L99: 14136.5 -> 14098.875
$ws.Cells.Item(99, 14).Value = -17094.875  # N99: -17132.5 -> -17094.875
$ws.Cells.Item(122, 8).Value = 1351  # H122: 1248 -> 1351
$ws.Cells.Item(122, 9).Value = 1329  # I122: 1214.3334 -> 1329
$ws.Cells.Item(122, 11).Value = 3987  # K122: 3643.0002 -> 3987
$ws.Cells.Item(122, 13).Value = -1537  # M122: -1193.0002 -> -1537
$ws.Cells.Item(126, 8).Value = 10799.517  # H126: 10818.936 -> 10799.517
$ws.Cells.Item(126, 10).Value = 14098.875  # J126: 14136.5 -> 14098.875
$ws.Cells.Item(126, 12).Value = 42296.625  # L126: 42409.5 -> 42296.625
$ws.Cells.Item(126, 14).Value = -47236.625  # N126: -47349.5 -> -47236.625
$ws.Cells.Item(134, 8).Value = 3997.5  # H134: 3353.6667 -> 3997.5
$ws.Cells.Item(134, 9).Value = 3996.3333  # I134: 2837.8 -> 3996.3333
$ws.Cells.Item(134, 10).Value = 3998.2  # J134: 3998.5 -> 3998.2
$ws.Cells.Item(134, 11).Value = 11988.9999  # K134: 8513.400000000001 -> 11988.9999
$ws.Cells.Item(134, 12).Value = 11994.6  # L134: 11995.5 -> 11994.6
$ws.Cells.Item(134, 13).Value = -9453.999899999999  # M134: -5978.400000000001 -> -9453.999899999999
$ws.Cells.Item(134, 14).Value = -17064.6  # N134: -17065.5 -> -17064.6
# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(93, 8).Value = 2913.5  # H93: 2314.2 -> 2913.5
$ws.Cells.Item(93, 9).Value = 0  # I93: 1744 -> 0
$ws.Cells.Item(93, 10).Value = 2913.5  # J93: 2456.75 -> 2913.5
$ws.Cells.Item(93, 11).Value = 0  # K93: 5232 -> 0
$ws.Cells.Item(93, 12).Value = 8740.5  # L93: 7370.25 -> 8740.5
$ws.Cells.Item(93, 13).ClearContents()  # M93: was -3360
$ws.Cells.Item(93, 14).Value = -12484.5  # N93: -11114.25 -> -12484.5
$ws.Cells.Item(122, 8).Value = 282.7143  # H122: 530 -> 282.7143
$ws.Cells.Item(122, 9).Value = 276.5  # I122: 280.66666 -> 276.5
$ws.Cells.Item(122, 10).Value = 291  # J122: 717 -> 291
$ws.Cells.Item(122, 11).Value = 2488.5  # K122: 2525.99994 -> 2488.5
$ws.Cells.Item(122, 12).Value = 2619  # L122: 6453 -> 2619
$ws.Cells.Item(122, 13).Value = -38.5  # M122: -75.9999399999997 -> -38.5
$ws.Cells.Item(122, 14).Value = -7519  # N122: -11353 -> -7519
$ws.Cells.Item(132, 8).Value = 3897.625  # H132: 3530.1667 -> 3897.625
$ws.Cells.Item(132, 10).Value = 5000  # J132: 0 -> 5000
$ws.Cells.Item(132, 12).Value = 45000  # L132: 0 -> 45000
$ws.Cells.Item(132, 14).Value = -50060  # N132: None -> -50060
# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(14, 8).Value = 427.2  # H14: 386.16666 -> 427.2
$ws.Cells.Item(14, 9).Value = 186  # I14: 183.5 -> 186
$ws.Cells.Item(14, 11).Value = 186  # K14: 183.5 -> 186
$ws.Cells.Item(14, 13).Value = -18  # M14: -15.5 -> -18
$ws.Cells.Item(20, 8).Value = 27479.062  # H20: 44444 -> 27479.062
$ws.Cells.Item(20, 9).Value = 9005  # I20: 0 -> 9005
$ws.Cells.Item(20, 10).Value = 28710.666  # J20: 44444 -> 28710.666
$ws.Cells.Item(20, 11).Value = 9005  # K20: 0 -> 9005
$ws.Cells.Item(20, 12).Value = 28710.666  # L20: 44444 -> 28710.666
$ws.Cells.Item(20, 13).Value = -8760  # M20: None -> -8760
$ws.Cells.Item(20, 14).Value = -29200.666  # N20: -44934 -> -29200.666
$ws.Cells.Item(35, 8).Value = 1682986.6  # H35: 2015984 -> 1682986.6
$ws.Cells.Item(35, 9).Value = 1682986.6  # I35: 2015984 -> 1682986.6
$ws.Cells.Item(35, 11).Value = 1682986.6  # K35: 2015984 -> 1682986.6
$ws.Cells.Item(35, 13).Value = -1682688.6  # M35: -2015686 -> -1682688.6
$ws.Cells.Item(102, 8).Value = 0  # H102: 1874 -> 0
$ws.Cells.Item(102, 10).Value = 0  # J102: 1874 -> 0
$ws.Cells.Item(102, 12).Value = 0  # L102: 1874 -> 0
$ws.Cells.Item(102, 14).ClearContents()  # N102: was -5118
$ws.Cells.Item(124, 8).Value = 0  # H124: 74959.5 -> 0
$ws.Cells.Item(124, 10).Value = 0  # J124: 74959.5 -> 0
$ws.Cells.Item(124, 12).Value = 0  # L124: 74959.5 -> 0
$ws.Cells.Item(124, 14).ClearContents()  # N124: was -84779.5
$ws.Cells.Item(126, 8).Value = 4494.5  # H126: 3333.3333 -> 4494.5
$ws.Cells.Item(126, 9).Value = 4494.5  # I126: 3500 -> 4494.5
$ws.Cells.Item(126, 10).Value = 0  # J126: 3000 -> 0
$ws.Cells.Item(126, 11).Value = 13483.5  # K126: 10500 -> 13483.5
$ws.Cells.Item(126, 12).Value = 0  # L126: 9000 -> 0
$ws.Cells.Item(126, 13).Value = -11013.5  # M126: -8030 -> -11013.5
$ws.Cells.Item(126, 14).ClearContents()  # N126: was -13940
# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 495  # H46: 446 -> 495
$ws.Cells.Item(46, 10).Value = 495  # J46: 446 -> 495
$ws.Cells.Item(46, 12).Value = 495  # L46: 446 -> 495
$ws.Cells.Item(46, 14).Value = -871  # N46: -822 -> -871
$ws.Cells.Item(61, 8).Value = 13892151  # H61: 12348990 -> 13892151
$ws.Cells.Item(61, 9).Value = 18521818  # I61: 15876202 -> 18521818
$ws.Cells.Item(61, 10).Value = 3149.5  # J61: 3751.5 -> 3149.5
$ws.Cells.Item(61, 11).Value = 18521818  # K61: 15876202 -> 18521818
$ws.Cells.Item(61, 12).Value = 3149.5  # L61: 3751.5 -> 3149.5
$ws.Cells.Item(61, 13).Value = -18521616  # M61: -15876000 -> -18521616
$ws.Cells.Item(61, 14).Value = -3553.5  # N61: -4155.5 -> -3553.5
$ws.Cells.Item(68, 8).Value = 5019.4  # H68: 4666.3335 -> 5019.4
$ws.Cells.Item(68, 9).Value = 3465.6667  # I68: 3324.5 -> 3465.6667
$ws.Cells.Item(68, 11).Value = 3465.6667  # K68: 3324.5 -> 3465.6667
$ws.Cells.Item(68, 13).Value = -2716.6667  # M68: -2575.5 -> -2716.6667
$ws.Cells.Item(71, 8).Value = 5019.4  # H71: 4666.3335 -> 5019.4
$ws.Cells.Item(71, 9).Value = 3465.6667  # I71: 3324.5 -> 3465.6667
$ws.Cells.Item(71, 11).Value = 17328.3335  # K71: 16622.5 -> 17328.3335
$ws.Cells.Item(71, 13).Value = -13584.3335  # M71: -12878.5 -> -13584.3335
$ws.Cells.Item(113, 8).Value = 13892151  # H113: 12348990 -> 13892151
$ws.Cells.Item(113, 9).Value = 18521818  # I113: 15876202 -> 18521818
$ws.Cells.Item(113, 10).Value = 3149.5  # J113: 3751.5 -> 3149.5
$ws.Cells.Item(113, 11).Value = 18521818  # K113: 15876202 -> 18521818
$ws.Cells.Item(113, 12).Value = 3149.5  # L113: 3751.5 -> 3149.5
$ws.Cells.Item(113, 13).Value = -18519648  # M113: -15874032 -> -18519648
$ws.Cells.Item(113, 14).Value = -7489.5  # N113: -8091.5 -> -7489.5
$ws.Cells.Item(122, 8).Value = 3990  # H122: 3995.6667 -> 3990
$ws.Cells.Item(122, 9).Value = 3990  # I122: 3997 -> 3990
$ws.Cells.Item(122, 10).Value = 0  # J122: 3993 -> 0
$ws.Cells.Item(122, 11).Value = 11970  # K122: 11991 -> 11970
$ws.Cells.Item(122, 12).Value = 0  # L122: 11979 -> 0
$ws.Cells.Item(122, 13).Value = -9520  # M122: -9541 -> -9520
$ws.Cells.Item(122, 14).ClearContents()  # N122: was -16879
$ws.Cells.Item(136, 8).Value = 7068.2  # H136: 7298.1113 -> 7068.2
$ws.Cells.Item(136, 9).Value = 6202.2  # I136: 6503 -> 6202.2
$ws.Cells.Item(136, 11).Value = 18606.6  # K136: 19509 -> 18606.6
$ws.Cells.Item(136, 13).Value = -16056.6  # M136: -16959 -> -16056.6
# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(3, 8).Value = 62900.6  # H3: 30222.555 -> 62900.6
$ws.Cells.Item(3, 9).Value = 130000  # I3: 52921.2 -> 130000
$ws.Cells.Item(3, 10).Value = 18167.666  # J3: 1849.25 -> 18167.666
$ws.Cells.Item(3, 11).Value = 130000  # K3: 52921.2 -> 130000
$ws.Cells.Item(3, 12).Value = 18167.666  # L3: 1849.25 -> 18167.666
$ws.Cells.Item(3, 13).Value = -129886  # M3: -52807.2 -> -129886
$ws.Cells.Item(3, 14).Value = -18395.666  # N3: -2077.25 -> -18395.666
$ws.Cells.Item(62, 8).Value = 4591.7856  # H62: 4618.933 -> 4591.7856
$ws.Cells.Item(62, 10).Value = 4999.2  # J62: 4999.1816 -> 4999.2
$ws.Cells.Item(62, 12).Value = 4999.2  # L62: 4999.1816 -> 4999.2
$ws.Cells.Item(62, 14).Value = -6247.2  # N62: -6247.1816 -> -6247.2
$ws.Cells.Item(65, 8).Value = 4591.7856  # H65: 4618.933 -> 4591.7856
$ws.Cells.Item(65, 10).Value = 4999.2  # J65: 4999.1816 -> 4999.2
$ws.Cells.Item(65, 12).Value = 24996  # L65: 24995.908 -> 24996
$ws.Cells.Item(65, 14).Value = -31236  # N65: -31235.908 -> -31236
$ws.Cells.Item(107, 8).Value = 423  # H107: 378.4 -> 423
$ws.Cells.Item(107, 9).Value = 500.66666  # I107: 425.5 -> 500.66666
$ws.Cells.Item(107, 11).Value = 1501.99998  # K107: 1276.5 -> 1501.99998
$ws.Cells.Item(107, 13).Value = 418.0000199999999  # M107: 643.5 -> 418.0000199999999
$ws.Cells.Item(122, 8).Value = 4022.1765  # H122: 4110.4707 -> 4022.1765
$ws.Cells.Item(122, 9).Value = 4089.3635  # I122: 4187.9 -> 4089.3635
$ws.Cells.Item(122, 10).Value = 3899  # J122: 3999.8572 -> 3899
$ws.Cells.Item(122, 11).Value = 12268.0905  # K122: 12563.7 -> 12268.0905
$ws.Cells.Item(122, 12).Value = 11697  # L122: 11999.5716 -> 11697
$ws.Cells.Item(122, 13).Value = -9818.0905  # M122: -10113.7 -> -9818.0905
$ws.Cells.Item(122, 14).Value = -16597  # N122: -16899.5716 -> -16597
$ws.Cells.Item(132, 8).Value = 2094.6365  # H132: 2215.4443 -> 2094.6365
$ws.Cells.Item(132, 9).Value = 1856.4445  # I132: 1943.7142 -> 1856.4445
$ws.Cells.Item(132, 11).Value = 5569.333500000001  # K132: 5831.142599999999 -> 5569.333500000001
$ws.Cells.Item(132, 13).Value = -3039.333500000001  # M132: -3301.142599999999 -> -3039.333500000001
